$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raids")

# Row 3: "The Ice Queens Reign" raid - update the story text (column B)
$ws.Range("B3").Value = "She awoke in the middle of the night. The snow fell upon the ground as her eyes adjusted to the darkness. She could see that she was upside down. The SUV her and her husband were driving had lost control on the snow filled roads. She looked over to see her husband, also upside down – alas he was dead. That’s when it hit her, She was alone – completely and utterly alone. Her son committed suicide a year before, and now her husband was dead. “Hello there” comes a voice. She looks over to see feet, a moment later she is standing beside the man, wearing a fedora. “I am The Poet, You must be Isabella” the man states. “My son, Have you seen my son?” She asks. It was all she could mutter as the snow fell in the darkness of the night."

# Row 3: raid_monster_ids (column D) - remove stray extra space before the comma
$ws.Range("D3").Value = "Corrupted Ice Mage,Queens Knight of the Ice Rose,Lost Memory of Her Son,Haunting Snowman,Living Ice sickle,Rabid Reindeer,Frozen King Krampus"

# The text edits above change the "best fit" column widths Excel stores for
# columns B (story) and D (raid_monster_ids). Set them to the recalculated
# widths explicitly (closest representable value for this engine's width
# rounding).
$ws.Columns.Item(2).ColumnWidth = 868.5
$ws.Columns.Item(4).ColumnWidth = 166.5
